$d = $word.ActiveDocument

function Replace-ParagraphXml($paragraphIndex, $xml) {
    $p = $d.Paragraphs($paragraphIndex)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $r.InsertXML($xml)
}

# 1) Delete trailing section (old paragraphs 19-23: border para, "3. Run compilation
#    script" heading, instructions paragraph, code block, and final paragraph) and
#    fold the closing sentence into paragraph 18 ("Wait while the files are created.")
$p18 = $d.Paragraphs(18)
$p23 = $d.Paragraphs($d.Paragraphs.Count)
$tailRange = $d.Range($p18.Range.Start, $p23.Range.End)
$tailRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Wait while the files</w:t></w:r><w:r><w:t xml:space="preserve"> are created.</w:t></w:r><w:r><w:t xml:space="preserve"> When it is done you will have a .csv file in the same directory as the input file. The output file’s name will be prepended with the date/time and appended with ‘*_compiled’.</w:t></w:r></w:p>')

# 2) Code block: workon dkifa_fslcalcs_env -> workon batch_fslstats_env (plus spell tags)
Replace-ParagraphXml 15 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:pBdr><w:top w:val="single" w:sz="6" w:space="6" w:color="E8E8E8"/><w:left w:val="single" w:sz="6" w:space="9" w:color="E8E8E8"/><w:bottom w:val="single" w:sz="6" w:space="6" w:color="E8E8E8"/><w:right w:val="single" w:sz="6" w:space="9" w:color="E8E8E8"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="EEEEFF"/><w:wordWrap w:val="0"/><w:ind w:left="180"/><w:rPr><w:rStyle w:val="HTMLCode"/><w:rFonts w:ascii="Consolas" w:eastAsiaTheme="majorEastAsia" w:hAnsi="Consolas"/><w:color w:val="333333"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="HTMLCode"/><w:rFonts w:ascii="Consolas" w:eastAsiaTheme="majorEastAsia" w:hAnsi="Consolas"/><w:color w:val="333333"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>workon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="HTMLCode"/><w:rFonts w:ascii="Consolas" w:eastAsiaTheme="majorEastAsia" w:hAnsi="Consolas"/><w:color w:val="333333"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="HTMLCode"/><w:rFonts w:ascii="Consolas" w:eastAsiaTheme="majorEastAsia" w:hAnsi="Consolas"/><w:color w:val="333333"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>batch_fslstats</w:t></w:r><w:r><w:rPr><w:rStyle w:val="HTMLCode"/><w:rFonts w:ascii="Consolas" w:eastAsiaTheme="majorEastAsia" w:hAnsi="Consolas"/><w:color w:val="333333"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>_env</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

# 3) Code block: bash run_fslstats_onall_dki_kfa.sh -> python compile_fsl_data.py
Replace-ParagraphXml 16 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:pBdr><w:top w:val="single" w:sz="6" w:space="6" w:color="E8E8E8"/><w:left w:val="single" w:sz="6" w:space="9" w:color="E8E8E8"/><w:bottom w:val="single" w:sz="6" w:space="6" w:color="E8E8E8"/><w:right w:val="single" w:sz="6" w:space="9" w:color="E8E8E8"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="EEEEFF"/><w:wordWrap w:val="0"/><w:spacing w:after="225"/><w:ind w:left="180"/><w:rPr><w:rStyle w:val="HTMLCode"/><w:rFonts w:ascii="Consolas" w:eastAsiaTheme="majorEastAsia" w:hAnsi="Consolas"/><w:color w:val="333333"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="HTMLCode"/><w:rFonts w:ascii="Consolas" w:eastAsiaTheme="majorEastAsia" w:hAnsi="Consolas"/><w:color w:val="333333"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>python compile_fsl_data.py</w:t></w:r></w:p>'

# 4) Note paragraph: dir() spell tag + .xlsx -> .csv
Replace-ParagraphXml 11 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Note: Its easiest to make this in MATLAB using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dir</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>() command and then saving the result as an .</w:t></w:r><w:r><w:t>csv</w:t></w:r><w:r><w:t xml:space="preserve"> table.</w:t></w:r></w:p>'

# 5) csv paragraph: input_file spell tag
Replace-ParagraphXml 10 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>First you need to put together a list of files</w:t></w:r><w:r><w:t xml:space="preserve"> to run FSL stats on. This must be a csv where the first row says </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>input_file</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and the remaining rows are full file paths to .nii files. Each .nii file will have its average value calculated.</w:t></w:r></w:p>'

# 6) Heading: "1. Set up a list of files to run fslstats on" -> spell tag around fslstats
Replace-ParagraphXml 9 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">1. Set up a list of files to run </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fslstats</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> on</w:t></w:r></w:p>'

# 7) "Before you run..." paragraph: merge trailing-space run into previous run
Replace-ParagraphXml 7 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Before you run, it might be a good idea to pull any new changes from the repo using git </w:t></w:r></w:p>'

# 8) "Use these directions..." paragraph: reworded
Replace-ParagraphXml 5 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Use these directions to get </w:t></w:r><w:r><w:t xml:space="preserve">mean values for a set of </w:t></w:r><w:r><w:t xml:space="preserve">.nii </w:t></w:r><w:r><w:t>images.</w:t></w:r><w:r><w:t xml:space="preserve"> You must first have done the setup.</w:t></w:r></w:p>'

# 9) Title: "DKI FSL Calcs: How to run scripts" -> "batch_fslstats: How to run scripts"
Replace-ParagraphXml 1 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>batch_fslstats</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: How to run scripts</w:t></w:r></w:p>'
